$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report table gained three new medicine rows. They need to be inserted
# right before the existing "totals" row (old row 6) and "footer" row
# (old row 7), which will shift down to rows 9 and 10 respectively.
$ws.Rows("6:8").Insert()

# Copy the formatting (styles) of the row-5 data row into the three new rows
# so they pick up the same cell styles (s="6".."9") as the existing data rows.
$ws.Range("A5:N5").Copy()
$ws.Range("A6:N8").PasteSpecial(-4122)  # xlPasteFormats

# Restore the row heights matching the pattern used by rows 4 (24.75) and 5/
# (25.5) - row 6 mirrors row 4's height, rows 7 & 8 mirror row 5's height.
$ws.Rows("6").RowHeight = 24.75
$ws.Rows("7").RowHeight = 25.5
$ws.Rows("8").RowHeight = 25.5

# Re-create the merged cell ranges for the new rows, same layout as rows 4/5.
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()

$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()

$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

# Fill in the three new medicine rows.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "MOTINORM 5 MG/5ML SUSP. 125 ML"
$ws.Range("H6").Value = "1:0"
$ws.Range("L6").Value = 31
$ws.Range("N6").Value = "1:0"

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "NANAZOXID 100MG/5ML PD. FOR ORAL SUSP. 60 ML"
$ws.Range("H7").Value = "2:0"
$ws.Range("L7").Value = 39
$ws.Range("N7").Value = "1:0"

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "ONDALENZ 4 MG 5 ORODISPERSIBLE FILMS"
$ws.Range("H8").Value = "0:3"
$ws.Range("L8").Value = 31.36
$ws.Range("N8").Value = "0:0"

# Update the totals row (now row 9) to reflect the new sum of all L column
# values (13 + 82 + 31 + 39 + 31.36 = 196.36).
$ws.Range("K9").Value = 196.36
